$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "요즘 Global MBA 시장 트렌드"
$ws.Range("E9").Value = "https://pdsi.pabii.com/global-mba-trend-2022/#utm_source=rss&utm_medium=rss&utm_campaign=global-mba-trend-2022"

$ws.Range("D44").Value = "principle (책) 리뷰 - Ray Dalio"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/123"

$ws.Range("D51").Value = "[css] input의 글자색, 배경색 바꾸기 (placeholder 포함)"
$ws.Range("E51").Value = "https://bskyvision.com/entry/css-input%EC%9D%98-%EA%B8%80%EC%9E%90%EC%83%89-%EB%B0%B0%EA%B2%BD%EC%83%89-%EB%B0%94%EA%BE%B8%EA%B8%B0-placeholder-%ED%8F%AC%ED%95%A8"

$wb.Save()
